$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.351.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.64%  '

# Row 3: 'Ethereum'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.51%  '

# Row 4: 'TetherUSD'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5: 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.34%  '

# Row 6: 'USDC'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '

# Row 7: 'XRP'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4452'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.00%  '

# Row 8: 'Cardano'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3720'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.04%  '

# Row 9: 'OKB'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.99'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.83%  '

# Row 10: 'Dogecoin'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07753'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.45%  '

# Row 11: 'Polygon'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.136'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.06%  '

# Row 12: 'Solana'
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '

# Row 13: 'BinanceUSD'
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.19'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.12%  '

# Row 14: 'Polkadot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.330'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.21%  '

# Row 15: 'Chainlink'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.629'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.24%  '

# Row 16: 'WrappedEther'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.855.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.20%  '

# Row 17: 'Litecoin'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.91%  '

# Row 18: 'ShibaInu'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001087'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.37%  '

# Row 19: 'TRON'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06544'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.64%  '

# Row 20: 'Dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9997'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '

# Row 21: 'Avalanche'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.60%  '

# Row 22: 'Uniswap'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.262'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.75%  '

# Row 23: 'WrappedBTC'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.396.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.67%  '

# Row 24: 'Cosmos'
$ws.Range('E24').Value = '  +3.14%  '

# Row 25: 'Toncoin'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.197'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.82%  '

# Row 26: 'EthereumClassic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.50%  '

# Row 27: 'Monero'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.02%  '

# Row 28: 'WrappedliquidstakedEther2.0'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.046.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.12%  '

# Row 29: 'LidoDAOToken'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.326'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.92%  '

# Row 30: 'BitcoinCash'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.01%  '

# Row 31: 'ImmutableX'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.214'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.66%  '

# Row 32: 'Filecoin'
$ws.Range('E32').Value = '  +5.15%  '

# Row 33: 'Stellar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09260'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.69%  '

# Row 34: 'HuobiToken'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.638'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.98%  '

# Row 35: 'Aptos'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.26'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.78%  '

# Row 36: 'VeChain'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02365'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.94%  '

# Row 37: 'Algorand'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2202'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.79%  '

# Row 38: 'InternetComputer(DFINITY)'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.216'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.79%  '

# Row 39: 'TheSandbox'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6631'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.99%  '

# Row 40: 'Hedera'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06252'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.44%  '

# Row 41: 'TrustWalletToken'
$ws.Range('E41').Value = '  +1.30%  '

# Row 42: 'FraxShare'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.198'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.84%  '

# Row 43: 'WEMIXTOKEN'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.438'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.20%  '

# Row 44: 'Frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '

# Row 45: 'EnergySwap'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.70%  '

# Row 46: 'Decentraland'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6171'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.01%  '

# Row 47: 'PancakeSwap'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.780'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.75%  '

# Row 48: 'NEARProtocol'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.051'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.83%  '

# Row 49: 'Quant'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '127.38'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.78%  '

# Row 50: 'EOS'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.162'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.84%  '

# Row 51: 'Cronos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07014'
$ws.Range('D51').Style = 'Normal'
